$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the first test case (row 2) with the new "room by guest capacity"
# test case (EPICA3. HU6 content), per the commit's sharedStrings/worksheet
# diff. The "EPICA /HU" (H2) and "Estado" (G2) columns are left untouched.
$nl = [char]10

$ws.Range("B2").Value2 = "Selección de habitacion por capacid"
$ws.Range("C2").Value2 = "Verificar la consulta de habitacion por cantidad maxima de huespedes"
$ws.Range("D2").Value2 = "1. Navegar a la página de consulta de habitaciones disponibles" + $nl + "2. Ingresar cantidad de huespedes" + $nl + "3. Realizar la consulta botón consultar."
$ws.Range("E2").Value2 = "Se muestran las habitaciones que tengan la capacidad especificada de huespedes"
$ws.Range("F2").Value2 = "Se muestran las habitaciones que tengan la capacidad especificada de huespedes"

# Match the author's final selection/scroll state captured in the saved file.
$ws.Activate()
$ws.Range("F2").Select()
